# EJ46FF data added and compiled
# Fill in the previously-empty tension table (columns B:G, rows 2:8) on the
# "EJ46FF" worksheet, then leave the selection on B3 (matches the saved
# workbook view in the target file).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EJ46FF")

$data = @{
    2 = @(286.4, 218.3, 170.7, 130.2, 94.8, 73.2)
    3 = @(292.2, 222.3, 174.4, 132.5, 96.8, 74.6)
    4 = @(298.2, 227.9, 179.2, 135.5, 99, 76.5)
    5 = @(303.4, 232.1, 183.9, 137.9, 101.2, 77.9)
    6 = @(309.7, 236.3, 187.7, 140.8, 103.4, 79.3)
    7 = @(316, 241.1, 191.4, 143.8, 105.4, 80.8)
    8 = @(322.4, 246.3, 195, 146.7, 107.5, 82.6)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}

$ws.Range("B3").Select()
